$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the B66 log entry text (finish wording for CPU chapter)
$ws.Range("B66").Value = "Added images of Ivy Bridge and Kepler architecture. Write chapter about CPU hardware architecture"

# Add the new log entry for row 67: date + description
$ws.Range("A67").Value = 41378
$ws.Range("B67").Value = "Finished chapter about the GPU hardware architecture. Added an additional paper about GPU optimizations"

# Move the active selection to the newly added entry
$ws.Range("B67").Select()
